$wb = $excel.ActiveWorkbook

# Insert a new "is_targeted list" sheet right after "analyte_class list"
# (and before "ms_source list"), matching the other "* list" helper sheets.
$afterSheet = $wb.Worksheets.Item("analyte_class list")
$listSheet = $wb.Worksheets.Add($null, $afterSheet)
$listSheet.Name = "is_targeted list"

# Populate it with the two allowed boolean-looking values as literal text
# (leading apostrophe forces text instead of Excel auto-coercing to a
# logical TRUE/FALSE value).
$listSheet.Range("A1").Value = "'TRUE"
$listSheet.Range("A2").Value = "'FALSE"

# Point the "is_targeted" column's validation at the new list sheet instead
# of the inline "TRUE,FALSE" literal formula, matching the other list-backed
# validations in this workbook.
$mainSheet = $wb.Worksheets.Item("Export as TSV")
$col = $mainSheet.Range("N2:N1048576")
$col.Validation.Modify(3, 1, 1, "='is_targeted list'!`$A`$1:`$A`$2")
$col.Validation.ErrorTitle = "Value must come from list"
$col.Validation.ErrorMessage = "Value must be one of: TRUE / FALSE."
